# The workbook tracks 轿车 (passenger car) sales/inventory stats, one row
# per "旬" (A/B/C/D) within each year. The B-旬 and C-旬 rows in every
# year-group were recorded in the wrong order; swap them back (A and D stay
# put). Columns F (轿车产销率) and G (轿车销售量) were retired entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-swap via Range.Copy (like a native Excel cut/paste) instead of Value2
# assignment: writing "" through Value2/Clear() actually removes the cell,
# while Copy()'ing a blank cell onto another cell leaves the destination's
# own pre-existing blank cell untouched (it does not get cleared out). Each
# group gets its own never-before-touched scratch row so that quirk can
# never carry stale data forward from a previous group.
$scratchBase = 200

$groupIndex = 0
for ($groupStart = 2; $groupStart -le 80; $groupStart += 4) {
    $bRow = $groupStart + 1
    $cRow = $groupStart + 2
    $scratchRow = $scratchBase + $groupIndex

    $ws.Range("A" + $bRow + ":E" + $bRow).Copy($ws.Range("A" + $scratchRow + ":E" + $scratchRow))
    $ws.Range("A" + $cRow + ":E" + $cRow).Copy($ws.Range("A" + $bRow + ":E" + $bRow))
    $ws.Range("A" + $scratchRow + ":E" + $scratchRow).Copy($ws.Range("A" + $cRow + ":E" + $cRow))

    $groupIndex++
}

$ws.Range("A" + $scratchBase + ":E" + ($scratchBase + $groupIndex - 1)).Clear()

# Drop the retired F/G columns (轿车产销率, 轿车销售量) entirely.
$ws.Columns.Item(6).Delete()
$ws.Columns.Item(6).Delete()
